$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 18,3
$arr[0,0] = 0.39397285841684893
$arr[0,1] = 0.34099433294687309
$arr[0,2] = -0.0044506466611495821
$arr[1,0] = 0.27831071534737412
$arr[1,1] = 0.31573146843522343
$arr[1,2] = -0.0044590127862502741
$arr[2,0] = 0.33185839385392346
$arr[2,1] = 0.37703067186178041
$arr[2,2] = -0.0045530038757724163
$arr[3,0] = 0.36291562613538619
$arr[3,1] = 0.35901250240432675
$arr[3,2] = -0.0045287913645162756
$arr[4,0] = 0.33185839385392346
$arr[4,1] = 0.37703067186178041
$arr[4,2] = -0.0045530038757724163
$arr[5,0] = 0.30508455460064876
$arr[5,1] = 0.34638107014850195
$arr[5,2] = -0.0046039436785600563
$arr[6,0] = 0.33185839385392346
$arr[6,1] = 0.37703067186178041
$arr[6,2] = -0.0045530038757724163
$arr[7,0] = 0.34069355018133618
$arr[7,1] = 0.36035918670473399
$arr[7,2] = -0.0045877960664201679
$arr[8,0] = 0.30508455460064876
$arr[8,1] = 0.34638107014850195
$arr[8,2] = -0.0046039436785600563
$arr[9,0] = 0.30495036946513049
$arr[9,1] = 0.30604904155629309
$arr[9,2] = -0.0045381668545978601
$arr[10,0] = 0.34069355018133618
$arr[10,1] = 0.36035918670473399
$arr[10,2] = -0.0045877960664201679
$arr[11,0] = 0.30508455460064876
$arr[11,1] = 0.34638107014850195
$arr[11,2] = -0.0046039436785600563
$arr[12,0] = 0.34069355018133618
$arr[12,1] = 0.36035918670473399
$arr[12,2] = -0.0045877960664201679
$arr[13,0] = 0.30508455460064876
$arr[13,1] = 0.34638107014850195
$arr[13,2] = -0.0046039436785600563
$arr[14,0] = 0.31391971092806148
$arr[14,1] = 0.32970958499145553
$arr[14,2] = -0.004610650116038912
$arr[15,0] = 0.30950213276435512
$arr[15,1] = 0.33804532756997874
$arr[15,2] = -0.0046094514063054429
$arr[16,0] = 0.31391971092806148
$arr[16,1] = 0.32970958499145553
$arr[16,2] = -0.004610650116038912
$arr[17,0] = 0.32730663055469883
$arr[17,1] = 0.34503438584809476
$arr[17,2] = -0.0046253997727040939

$ws.Range("A31:C48").Value = $arr

